$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The cell that held the single user-story paragraph:
#   "As a developer, I want to build a stock database of products so that
#    the sales manager can check all the available products, its quantity
#    and price at the same time in stock database"
# becomes two paragraphs:
#   1) "As a user, I want to check the stock database"
#   2) "of products so that I can check all the available products, its
#       quantity and price at the same time in stock database"
# with the run font size changed from 16pt (sz 32) to 14pt (sz 28).
# ---------------------------------------------------------------------------

$oldText = "As a developer, I want to build a stock database of products so that the sales manager can check all the available products, its quantity and price at the same time in stock database"

$newPara1 = "As a user, I want to check the stock database"
$newPara2 = "of products so that I can check all the available products, its quantity and price at the same time in stock database"

# Step 1: locate the paragraph and split its text into two paragraphs.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$oldText*") {
        $p.Range.Text = $newPara1 + "`r" + $newPara2
    }
}

# Step 2: re-locate paragraph 1 ("As a user...") and recreate the run
# boundaries seen in the target document by nudging the font size of each
# run span in turn (this both sets the new 14pt/sz-28 size and keeps the
# runs from being merged back together).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "$newPara1*") {
        $s = $p.Range.Start
        $d.Range($s + 0,  $s + 22).Font.Size = 14
        $d.Range($s + 22, $s + 45).Font.Size = 14
        break
    }
}

# Step 3: re-locate paragraph 2 ("of products...") and recreate its four
# run spans the same way.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "$newPara2*") {
        $s = $p.Range.Start
        $d.Range($s + 0,   $s + 1).Font.Size = 14
        $d.Range($s + 1,   $s + 11).Font.Size = 14
        $d.Range($s + 11,  $s + 20).Font.Size = 14
        $d.Range($s + 20,  $s + 117).Font.Size = 14
        break
    }
}

# Step 4: apply the matching complex-script size (szCs) across each whole
# paragraph in one shot -- setting SizeBi on a sub-range alone does not take
# effect in this host, but applying it across the already-split paragraph
# range updates every run inside it.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "$newPara1*") {
        $p.Range.Font.SizeBi = 14
    }
    if ($p.Range.Text -like "$newPara2*") {
        $p.Range.Font.SizeBi = 14
    }
}
